# Applies the weekly fruit/vegetable data reshuffle described by the diff.
# Columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) are updated
# for rows 3,4,5,6,7,8,9,12,13,14 with the "after" values taken from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = @{ D = 44446; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 467 }
    4  = @{ D = 44460; J = 45; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 433 }
    5  = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
    6  = @{ D = 44376; J = 25; K = 18000; L = 18000; M = 18000; O = "Provincia de Limarí"; P = 600 }
    7  = @{ D = 44453; J = 50; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
    8  = @{ D = 44432; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 467 }
    9  = @{ D = 44449; J = 45; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
    12 = @{ D = 44474; J = 45; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí"; P = 333 }
    13 = @{ D = 44418; J = 30; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 500 }
    14 = @{ D = 44421; J = 25; K = 15000; L = 16000; M = 15400; O = "Provincia de Limarí"; P = 513 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
}
